$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.020335
$ws.Range("H2").Value = 0.061005
$ws.Range("I2").Value = 0.009804808687698561
$ws.Range("J2").Value = 0.009804808687698559
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 8.081040666666667
$ws.Range("N2").Value = 24.243122
$ws.Range("O2").Value = 0.4661250698616886
$ws.Range("P2").Value = 0.4661250698616886
$ws.Range("Q2").Value = 0.1643279619566667
$ws.Range("R2").Value = 1.47895165761
$ws.Range("S2").Value = 0.004570267134533983
$ws.Range("T2").Value = 0.004570267134533982
$ws.Range("G3").Value = 0.020335
$ws.Range("H3").Value = 0.061005
$ws.Range("I3").Value = 0.009804808687698561
$ws.Range("J3").Value = 0.009804808687698559
$ws.Range("O3").Value = 0.4037865631294714
$ws.Range("P3").Value = 0.4037865631294715
$ws.Range("Q3").Value = 0.1423511140566667
$ws.Range("R3").Value = 1.28116002651
$ws.Range("S3").Value = 0.003959050002147785
$ws.Range("T3").Value = 0.003959050002147785
$ws.Range("G4").Value = 0.020335
$ws.Range("H4").Value = 0.061005
$ws.Range("I4").Value = 0.009804808687698561
$ws.Range("J4").Value = 0.009804808687698559
$ws.Range("O4").Value = 0.1300883670088399
$ws.Range("P4").Value = 0.1300883670088399
$ws.Range("Q4").Value = 0.04586141704666666
$ws.Range("R4").Value = 0.41275275342
$ws.Range("S4").Value = 0.001275491551016792
$ws.Range("T4").Value = 0.001275491551016792
$ws.Range("I5").Value = 0.1486140913768632
$ws.Range("J5").Value = 0.1486140913768632
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 8.081040666666667
$ws.Range("N5").Value = 24.243122
$ws.Range("O5").Value = 0.4661250698616886
$ws.Range("P5").Value = 0.4661250698616886
$ws.Range("Q5").Value = 2.490762597402
$ws.Range("R5").Value = 22.416863376618
$ws.Range("S5").Value = 0.06927275372547176
$ws.Range("T5").Value = 0.06927275372547174
$ws.Range("I6").Value = 0.1486140913768632
$ws.Range("J6").Value = 0.1486140913768632
$ws.Range("O6").Value = 0.4037865631294714
$ws.Range("P6").Value = 0.4037865631294715
$ws.Range("S6").Value = 0.06000837318967282
$ws.Range("T6").Value = 0.06000837318967282
$ws.Range("I7").Value = 0.1486140913768632
$ws.Range("J7").Value = 0.1486140913768632
$ws.Range("O7").Value = 0.1300883670088399
$ws.Range("P7").Value = 0.1300883670088399
$ws.Range("S7").Value = 0.01933296446171865
$ws.Range("T7").Value = 0.01933296446171865
$ws.Range("H8").Value = 5.236273000000001
$ws.Range("I8").Value = 0.8415810999354383
$ws.Range("J8").Value = 0.8415810999354382
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 8.081040666666667
$ws.Range("N8").Value = 24.243122
$ws.Range("O8").Value = 0.4661250698616886
$ws.Range("P8").Value = 0.4661250698616886
$ws.Range("Q8").Value = 14.10484501825622
$ws.Range("R8").Value = 126.943605164306
$ws.Range("S8").Value = 0.3922820490016829
$ws.Range("T8").Value = 0.3922820490016828
$ws.Range("H9").Value = 5.236273000000001
$ws.Range("I9").Value = 0.8415810999354383
$ws.Range("J9").Value = 0.8415810999354382
$ws.Range("O9").Value = 0.4037865631294714
$ws.Range("P9").Value = 0.4037865631294715
$ws.Range("S9").Value = 0.3398191399376508
$ws.Range("T9").Value = 0.3398191399376508
$ws.Range("H10").Value = 5.236273000000001
$ws.Range("I10").Value = 0.8415810999354383
$ws.Range("J10").Value = 0.8415810999354382
$ws.Range("O10").Value = 0.1300883670088399
$ws.Range("P10").Value = 0.1300883670088399
$ws.Range("Q10").Value = 3.936446190036889
$ws.Range("R10").Value = 35.42801571033201
$ws.Range("S10").Value = 0.1094799109961045
$ws.Range("T10").Value = 0.1094799109961045
